$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Insert a new column before column C (Chuyen Nganh / Major) on sheet2
$ws2.Columns("C").Insert()

$ws2.Range("C1").Value = "Chuyên Ngành"
$ws2.Range("C2").Value = "Công Nghệ Thông Tin"
$ws2.Range("C3").Value = "Công Nghệ Thông Tin"
$ws2.Range("C4").Value = "Quản Trị Kinh Doanh"
$ws2.Range("C5").Value = "Công Nghệ Thông Tin"
$ws2.Range("C6").Value = "Quản Trị Kinh Doanh"

$ws2.Columns("C").ColumnWidth = 23.3333333333

# Update selections: sheet2 keeps a selection on the new column,
# while sheet1 becomes the active (tab-selected) sheet.
$ws2.Activate()
$ws2.Range("C2").Select()

$ws1.Activate()
$ws1.Range("F9").Select()
